$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-23 Wednesday" "2024-10-24 Thursday"

Replace-Text "39÷2=19, 1" "99÷9=11, 0"
Replace-Text "28÷5=5, 3" "76÷6=12, 4"
Replace-Text "62÷3=20, 2" "42÷2=21, 0"
Replace-Text "86÷6=14, 2" "95÷5=19, 0"
Replace-Text "68÷4=17, 0" "38÷7=5, 3"

Replace-Text "73÷7=10, 3" "49÷8=6, 1"
Replace-Text "18÷6=3, 0" "70÷7=10, 0"
Replace-Text "47÷4=11, 3" "80÷5=16, 0"
Replace-Text "44÷2=22, 0" "60÷7=8, 4"
Replace-Text "25÷4=6, 1" "12÷4=3, 0"

Replace-Text "41÷4=10, 1" "31÷3=10, 1"
Replace-Text "77÷6=12, 5" "61÷3=20, 1"
Replace-Text "14÷7=2, 0" "83÷9=9, 2"
Replace-Text "68÷6=11, 2" "60÷8=7, 4"
Replace-Text "30÷5=6, 0" "98÷3=32, 2"

Replace-Text "20÷7=2, 6" "94÷2=47, 0"
Replace-Text "33÷6=5, 3" "18÷4=4, 2"
Replace-Text "60÷4=15, 0" "92÷8=11, 4"
Replace-Text "22÷2=11, 0" "62÷9=6, 8"
Replace-Text "79÷3=26, 1" "10÷8=1, 2"

Replace-Text "98÷9=10, 8" "26÷3=8, 2"
Replace-Text "53÷7=7, 4" "70÷6=11, 4"
Replace-Text "55÷2=27, 1" "15÷5=3, 0"
Replace-Text "86÷8=10, 6" "26÷4=6, 2"
Replace-Text "21÷3=7, 0" "97÷5=19, 2"
